# Update the EmailID in cell A3 on every monthly sheet (attendance roster)
# from the old address to the new one, per the commit's requested change.

$wb = $excel.ActiveWorkbook

$newEmail = "avinash18dce.kumar@gmail.com"

foreach ($ws in $wb.Worksheets) {
    $ws.Range("A3").Value = $newEmail
}
